# Covbot "queries_by_difficulty" worksheet update.
# Column A holds the subset of shared-string questions that have already
# been worked through (highlighted with the built-in "Good" style); three
# more rows (Serbia 2021 / Russia 2018 / peak-confirmed-cases) were marked
# done, which pushes the remaining "done" questions in A down by three rows
# and leaves the last three rows of that block blank again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column A: give rows 24-28 the same "Good" highlight as the rows above them ---
foreach ($r in 24..28) {
    $ws.Range("A$r").Style = "Good"
}

# --- Column A: newly answered questions now occupy rows 29-31 ---
$ws.Range("A29").Value = "How many new cases were discovered in Serbia in 2021?"
$ws.Range("A30").Value = "How many new cases were detected in Russia in 2018?"
$ws.Range("A31").Value = "What is the peak number of confirmed cases in certain country"
foreach ($r in 29..31) {
    $ws.Range("A$r").Style = "Good"
}

# --- Column A: everything that used to start at row 32 shifts up by three rows ---
$ws.Range("A32").Value = "What is the peak number of vaccinated people in a day of a certain country"
$ws.Range("A33").Value = "How many new cases of COVID are there today in Hong Kong?"
$ws.Range("A34").Value = "How many new cases of COVID are this week in Hong Kong?"
$ws.Range("A35").Value = "Which country has had the most corona cases?"
$ws.Range("A36").Value = "On which day were the most cases reported?"
$ws.Range("A37").Value = "In which country did most people get vaccinated?"
foreach ($r in 32..37) {
    $ws.Range("A$r").Style = "Good"
}

# --- Column A: rows 38-40 no longer have a "done" entry ---
foreach ($r in 38..40) {
    $ws.Range("A$r").ClearContents()
}

# --- View state: zoomed in on the newly-edited block instead of scrolled to A13 ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 125
$ws.Range("A33:A37").Select()
